$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date header for column C, same style as B1 (bold/bordered/centered)
$ws.Range("C1").Value = "13-01-2023"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Reorder rows 2-8: funds first (alphabetical-ish as in target), then avg, then total
# and add the new "13-01-2023" column (C) values alongside existing (B) values.
$ws.Range("A2").Value = "Delta Acciones"
$ws.Range("B2").Value = 63857.43
$ws.Range("C2").Value = 63891.57

$ws.Range("A3").Value = "Delta Select"
$ws.Range("B3").Value = 82239.83
$ws.Range("C3").Value = 83732.59

$ws.Range("A4").Value = "Delta gestion V"
$ws.Range("B4").Value = 29533.04
$ws.Range("C4").Value = 29240.39

$ws.Range("A5").Value = "Fima Acciones"
$ws.Range("B5").Value = 1940.65
$ws.Range("C5").Value = 2199.17

$ws.Range("A6").Value = "Fima PB Acciones"
$ws.Range("B6").Value = 3895.64
$ws.Range("C6").Value = 4334.54

$ws.Range("A7").Value = "avg"
$ws.Range("B7").Value = 36293.32
$ws.Range("C7").Value = 36679.65

$ws.Range("A8").Value = "total"
$ws.Range("B8").Value = 181466.59
$ws.Range("C8").Value = 183398.26
